$wb = $excel.ActiveWorkbook

# --- Update the selection left on "Piper" first (must happen while it is
#     still the active sheet, before the new sheet steals focus) -----------
$piper = $wb.Worksheets.Item("Piper")
$piper.Select()
$piper.Range("A1:F1").Select()

# --- Insert the new "Ipanema" worksheet right before "Piper" --------------
$ws = $wb.Worksheets.Add($piper)
$ws.Name = "Ipanema"

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "Model"
$ws.Range("A1").NumberFormat = "@"
$ws.Range("B1").Value = "Capacity (US GAL)"
$ws.Range("C1").Value = "Wingspan (FT)"
$ws.Range("D1").Value = "70% boom (FT)"
$ws.Range("E1").Value = "Ideal height (FT)"
$ws.Range("F1").Value = "Working speed (MPH)"

# --- Data rows: Embraer EMB 201 / 201A / 202 / 202A / 203 -----------------
$ws.Range("A2").Value = 201
$ws.Range("B2").Value = 180
$ws.Range("C2").Value = 36.3

$ws.Range("A3").Value = "201A"
$ws.Range("B3").Value = 180
$ws.Range("C3").Value = 36.3

$ws.Range("A4").Value = 202
$ws.Range("B4").Value = 250
$ws.Range("C4").Value = 38.4

$ws.Range("A5").Value = "202A"
$ws.Range("B5").Value = 250
$ws.Range("C5").Value = 38.4

$ws.Range("A6").Value = 203
$ws.Range("B6").Value = 277
$ws.Range("C6").Value = 43.6

# --- Column widths (approximate best-fit, matches the saved workbook) -----
$ws.Columns.Item(1).ColumnWidth = 5.4987
$ws.Columns.Item(2).ColumnWidth = 15.1667
$ws.Columns.Item(3).ColumnWidth = 12.3307
$ws.Columns.Item(4).ColumnWidth = 12.8307
$ws.Columns.Item(5).ColumnWidth = 13.8307
$ws.Columns.Item(6).ColumnWidth = 18.3307

# --- View state: the new sheet becomes the active / selected tab ----------
$ws.Activate()
$excel.ActiveWindow.Zoom = 99
$ws.Range("F8").Select()
